$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 1095
    $ws.Range("F5").Value = 4690
    $ws.Range("F8").Value = 1413
    $ws.Range("F11").Value = 1220
    $ws.Range("F13").Value = 682
}
